$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Simple single-cell updates (B column renumbering)
$ws.Range("B21").Value = 87163
$ws.Range("B22").Value = 87047

# Rows 23 and 24 swap most of their content (A, D, E, F, G, H, Q, R),
# while B gets distinct new values.
$ws.Range("A23").Value = 130807440
$ws.Range("B23").Value = 87198
$ws.Range("D23").Value = "VU"
$ws.Range("E23").Value = 6003296
$ws.Range("F23").Value = "Stor odörspindling"
$ws.Range("G23").Value = "Cortinarius mussivus"
$ws.Range("H23").Value = "(Fr.) Melot"
$ws.Range("Q23").Value = 704377
$ws.Range("R23").Value = 6361495

$ws.Range("A24").Value = 130807443
$ws.Range("B24").Value = 90538
$ws.Range("D24").Value = "NT"
$ws.Range("E24").Value = 970
$ws.Range("F24").Value = "Bittermusseron"
$ws.Range("G24").Value = "Leucopaxillus gentianeus"
$ws.Range("H24").Value = "(Quél.) Kotl."
$ws.Range("Q24").Value = 704276
$ws.Range("R24").Value = 6361505

$ws.Range("B25").Value = 87071
$ws.Range("B26").Value = 87223
$ws.Range("B27").Value = 87210
$ws.Range("B28").Value = 91221
$ws.Range("B29").Value = 87210
$ws.Range("B30").Value = 87094
$ws.Range("B31").Value = 87198
$ws.Range("B32").Value = 87094
